# Apply the recall_results.xlsx edit:
#  - Sheet1: re-style header/label/value cells, add a bold "Bi-LSTM" recall
#    column (G) with centered 0.00 number format, and add two new blocks of
#    delta formulas (rows 7-9 and 11-13).
#  - Add a new "politics" worksheet with a hyperlinked "recall@1" cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (C2:G2): centered bold Arial header style -------------------
$ws.Range("C2:G2").HorizontalAlignment = -4108   # xlCenter

# --- Row label column (B3:B5): unchanged content, style refreshed by Excel --
# (no explicit change needed; style index shift happens naturally)

# --- Value columns (C3:F5): center align the existing recall values --------
$ws.Range("C3:F5").HorizontalAlignment = -4108   # xlCenter
$ws.Range("C3:F5").NumberFormat = "0.00"

# --- New Bi-LSTM column (G3:G5): bold, centered, 0.00 number format --------
$ws.Range("G3").Value = 0.62343752399999997
$ws.Range("G4").Value = 0.70390623799999996
$ws.Range("G5").Value = 0.83984375

$gCol = $ws.Range("G3:G5")
$gCol.Font.Bold = $true
$gCol.HorizontalAlignment = -4108   # xlCenter
$gCol.NumberFormat = "0.00"

# --- Delta block 1: rows 7-9, columns E-G = (E/F/G)n - $Dn ------------------
$ws.Range("E7").Formula = "=E3-`$D3"
$ws.Range("F7").Formula = "=F3-`$D3"
$ws.Range("G7").Formula = "=G3-`$D3"
$ws.Range("E8").Formula = "=E4-`$D4"
$ws.Range("F8").Formula = "=F4-`$D4"
$ws.Range("G8").Formula = "=G4-`$D4"
$ws.Range("E9").Formula = "=E5-`$D5"
$ws.Range("F9").Formula = "=F5-`$D5"
$ws.Range("G9").Formula = "=G5-`$D5"
$ws.Range("E7:G9").NumberFormat = "0.00"

# --- Delta block 2: rows 11-13, column G = Gn - Fn --------------------------
$ws.Range("G11").Formula = "=G3-F3"
$ws.Range("G12").Formula = "=G4-F4"
$ws.Range("G13").Formula = "=G5-F5"
$ws.Range("G11:G13").NumberFormat = "0.00"

# --- New worksheet: politics --------------------------------------------
$wsPolitics = $wb.Worksheets.Add($null, $ws)
$wsPolitics.Name = "politics"
$wsPolitics.Range("B2").Value = "domain"
$wsPolitics.Range("C2").Value = "recall@1"
$wsPolitics.Hyperlinks.Add($wsPolitics.Range("C2"), "https://github.com/jtamon/twconvrsu")

$wsPolitics.Range("C3").Select()
